# Insert a new weekly price record as row 26 ("Hortaliza, Terminal La Palmera
# de La Serena - Sandia"), pushing the existing rows 26-54 down to 27-55.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("26").Insert()

$ws.Range("A26").Value = 8
$ws.Range("B26").Value = "Terminal La Palmera de La Serena"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44494
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100112028
$ws.Range("G26").Value = "Sandia"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 700
$ws.Range("L26").Value = 750
$ws.Range("M26").Value = 725
$ws.Range("N26").Value = "$/kilo (volumen en unidades)"
$ws.Range("O26").Value = "Perú"
$ws.Range("P26").Value = 725
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"
